# Connected the create timetable button to the script
# Updates the generated timetable cells on each day-sheet (mon/tue/wed/thur/fri)
# to reflect the latest output of the timetable-generation script.

$wb = $excel.ActiveWorkbook

# --- Monday ---
$ws = $wb.Worksheets.Item("mon")
$ws.Range("E4").Value = "CSC423"
$ws.Range("G4").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("I13").Value = "CSC442"
$ws.Range("J13").Value = "CSC442"
$ws.Range("G14").Value = ""
$ws.Range("B20").Value = "GST111"
$ws.Range("C20").Value = "GST111"
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""
$ws.Range("I21").Value = ""
$ws.Range("J21").Value = "MAT112"
$ws.Range("K21").Value = "MAT112"
$ws.Range("H26").Value = ""
$ws.Range("I26").Value = ""

# --- Tuesday ---
$ws = $wb.Worksheets.Item("tue")
$ws.Range("C10").Value = "BIO111"
$ws.Range("D10").Value = "BIO111"
$ws.Range("J14").Value = ""
$ws.Range("J15").Value = "CSC425"
$ws.Range("H20").Value = ""
$ws.Range("I20").Value = "MAT111"
$ws.Range("J20").Value = "MAT111"
$ws.Range("J21").Value = "CSC424"
$ws.Range("J25").Value = "CIS421"
$ws.Range("K25").Value = "CIS421"

# --- Wednesday ---
$ws = $wb.Worksheets.Item("wed")
$ws.Range("I4").Value = "CSC424"
$ws.Range("J4").Value = "CSC424"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("G19").Value = ""
$ws.Range("D20").Value = "CIT111"
$ws.Range("E20").Value = "CIT111"
$ws.Range("F20").Value = ""
$ws.Range("I20").Value = ""
$ws.Range("J20").Value = ""
$ws.Range("E21").Value = "MAT111"
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = ""
$ws.Range("G25").Value = "CSC111"
$ws.Range("H25").Value = "CSC111"
$ws.Range("I25").Value = ""
$ws.Range("J25").Value = ""

# --- Thursday ---
$ws = $wb.Worksheets.Item("thur")
$ws.Range("D7").Value = "CSC441"
$ws.Range("E7").Value = "CSC441"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("J11").Value = "CSC424"
$ws.Range("K11").Value = "CSC424"
$ws.Range("F16").Value = ""
$ws.Range("G16").Value = ""
$ws.Range("G17").Value = ""
$ws.Range("F20").Value = "CSC442"
$ws.Range("I20").Value = "CHM111"
$ws.Range("J20").Value = "CHM111"
$ws.Range("F21").Value = ""
$ws.Range("J21").Value = "EDS421"

# --- Friday ---
$ws = $wb.Worksheets.Item("fri")
$ws.Range("C20").Value = ""
$ws.Range("E20").Value = "TMC111"
$ws.Range("F20").Value = ""
$ws.Range("C21").Value = "TMC421"
